$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
  @(366, @(@(6,29),@(7,19),@(8,9),@(9,9),@(10,1),@(11,-7),@(12,82),@(13,31),@(14,5),@(15,1019),@(16,1015),@(17,1012),@(18,11),@(19,9),@(20,4),@(21,26),@(22,5),@(24,0),@(25,6),@(26,"Rain"),@(27,123))),
  @(367, @(@(6,29),@(7,23),@(8,18),@(9,16),@(10,5),@(11,-1),@(12,48),@(13,29),@(14,9),@(15,1013),@(16,1010),@(17,1007),@(18,10),@(19,7),@(20,6),@(21,32),@(22,16),@(23,53),@(24,0),@(25,6),@(26,"Rain-Thunderstorm"),@(27,154))),
  @(368, @(@(6,28),@(7,23),@(8,18),@(9,15),@(10,7),@(11,-1),@(12,83),@(13,37),@(14,9),@(15,1014),@(16,1011),@(17,1009),@(18,10),@(19,8),@(20,3),@(21,26),@(22,13),@(23,35),@(24,0),@(25,4),@(27,288))),
  @(369, @(@(6,26),@(7,21),@(8,16),@(9,7),@(10,3),@(11,-5),@(12,46),@(13,29),@(14,13),@(15,1019),@(16,1016),@(17,1014),@(18,10),@(19,10),@(20,10),@(21,19),@(22,11),@(24,0),@(25,1),@(27,299))),
  @(370, @(@(6,28),@(7,20),@(8,13),@(9,12),@(10,4),@(11,-1),@(12,63),@(13,34),@(14,11),@(15,1021),@(16,1019),@(17,1017),@(18,10),@(19,10),@(20,10),@(21,14),@(22,6),@(24,0),@(27,320))),
  @(371, @(@(6,32),@(7,22),@(8,12),@(9,11),@(10,6),@(11,0),@(12,82),@(13,37),@(14,13),@(15,1020),@(16,1018),@(17,1016),@(18,10),@(19,8),@(20,6),@(21,8),@(22,2),@(24,0),@(27,319))),
  @(372, @(@(6,28),@(7,22),@(8,16),@(9,10),@(10,3),@(11,-5),@(12,55),@(13,29),@(14,7),@(15,1019),@(16,1015),@(17,1012),@(18,10),@(19,10),@(20,8),@(21,19),@(22,5),@(24,0),@(25,8),@(27,84))),
  @(373, @(@(6,26),@(7,21),@(8,16),@(9,14),@(10,8),@(11,-6),@(12,88),@(13,50),@(14,12),@(15,1014),@(16,1011),@(17,1008),@(18,10),@(19,7),@(20,3),@(21,26),@(22,10),@(23,50),@(24,7.87),@(25,8),@(26,"Rain-Thunderstorm"),@(27,114))),
  @(374, @(@(6,23),@(7,19),@(8,15),@(9,19),@(10,17),@(11,12),@(12,88),@(13,77),@(14,50),@(15,1012),@(16,1009),@(17,1006),@(18,10),@(19,5),@(20,1),@(21,19),@(22,5),@(24,3.05),@(25,6),@(26,"Rain"),@(27,103))),
  @(375, @(@(6,25),@(7,21),@(8,17),@(9,17),@(10,11),@(11,5),@(12,83),@(13,54),@(14,25),@(15,1015),@(16,1012),@(17,1008),@(18,10),@(19,8),@(20,5),@(21,26),@(22,14),@(24,0),@(25,3),@(27,288))),
  @(376, @(@(6,26),@(7,19),@(8,13),@(9,9),@(10,6),@(11,0),@(12,72),@(13,40),@(14,12),@(15,1017),@(16,1015),@(17,1013),@(18,10),@(19,10),@(20,10),@(21,19),@(22,11),@(24,0),@(27,322))),
  @(377, @(@(6,26),@(7,19),@(8,13),@(9,8),@(10,4),@(11,-3),@(12,63),@(13,37),@(14,15),@(15,1015),@(16,1014),@(17,1011),@(18,11),@(19,10),@(20,10),@(21,23),@(22,13),@(24,0),@(25,4),@(27,292))),
  @(378, @(@(6,25),@(7,19),@(8,12),@(9,7),@(10,3),@(11,-1),@(12,59),@(13,35),@(14,10),@(15,1016),@(16,1014),@(17,1013),@(18,11),@(19,10),@(20,10),@(21,23),@(22,11),@(24,0),@(25,3),@(27,288))),
  @(379, @(@(6,26),@(7,18),@(8,10),@(9,8),@(10,5),@(11,-3),@(12,76),@(13,40),@(14,15),@(15,1019),@(16,1016),@(17,1015),@(18,10),@(19,9),@(20,8),@(21,19),@(22,3),@(24,0),@(27,321))),
  @(380, @(@(6,28),@(7,19),@(8,11),@(9,11),@(10,6),@(11,2),@(12,82),@(13,39),@(14,13),@(15,1019),@(16,1016),@(17,1014),@(18,10),@(19,10),@(20,10),@(21,11),@(22,2),@(24,0),@(27,324))),
  @(381, @(@(6,30),@(7,21),@(8,12),@(9,12),@(10,7),@(11,3),@(12,77),@(13,40),@(14,14),@(15,1017),@(16,1014),@(17,1012),@(18,10),@(19,8),@(20,6),@(21,8),@(22,2),@(24,0),@(27,31))),
  @(382, @(@(6,32),@(7,22),@(8,12),@(9,13),@(10,8),@(11,4),@(12,77),@(13,42),@(14,12),@(15,1013),@(16,1011),@(17,1009),@(18,10),@(19,10),@(20,10),@(21,8),@(22,2),@(24,0),@(27,323))),
  @(383, @(@(6,33),@(7,24),@(8,16),@(9,14),@(10,9),@(11,-2),@(12,77),@(13,42),@(14,11),@(15,1012),@(16,1010),@(17,1009),@(18,10),@(19,7),@(20,5),@(21,11),@(22,3),@(24,0),@(27,275))),
  @(384, @(@(6,33),@(7,24),@(8,16),@(9,18),@(10,12),@(11,7),@(12,68),@(13,42),@(14,15),@(15,1012),@(16,1010),@(17,1008),@(18,10),@(19,8),@(20,6),@(21,11),@(22,3),@(24,0),@(27,280))),
  @(385, @(@(6,33),@(7,26),@(8,18),@(9,15),@(10,12),@(11,7),@(12,73),@(13,38),@(14,18),@(15,1012),@(16,1010),@(17,1007),@(18,10),@(19,9),@(20,6),@(21,19),@(22,6),@(24,0),@(25,1),@(27,285))),
  @(386, @(@(6,34),@(7,26),@(8,18),@(9,14),@(10,9),@(11,3),@(12,73),@(13,32),@(14,13),@(15,1011),@(16,1009),@(17,1008),@(18,11),@(19,10),@(20,10),@(21,19),@(22,8),@(24,0),@(27,295))),
  @(387, @(@(6,36),@(7,27),@(8,17),@(9,14),@(10,10),@(11,7),@(12,53),@(13,34),@(14,12),@(15,1013),@(16,1011),@(17,1010),@(18,10),@(19,9),@(20,7),@(21,14),@(22,6),@(24,0),@(25,4),@(27,170))),
  @(388, @(@(6,30),@(7,26),@(8,22),@(9,14),@(10,10),@(11,6),@(12,50),@(13,35),@(14,17),@(15,1014),@(16,1011),@(17,1008),@(18,10),@(19,7),@(20,5),@(21,23),@(22,6),@(23,32),@(24,0),@(25,6),@(26,"Rain"),@(27,106)))
)

foreach ($rowEntry in $rows) {
  $r = $rowEntry[0]
  $cells = $rowEntry[1]
  foreach ($cellEntry in $cells) {
    $c = $cellEntry[0]
    $v = $cellEntry[1]
    $ws.Cells.Item($r, $c).Value = $v
  }
}

[void]$ws.Range("F1").Select()

Write-Host "done"